# Implemented getting number of lines for methods and classes.
#
# 1) Re-order the rows of the "classFields" sheet to reflect the refreshed
#    field-declaration order coming out of the analysis tool (same set of
#    field rows, new relative order within each class).
# 2) Add two new sheets - "classNumberOfLines" and "methodNumberOfLines" -
#    holding the line-count metrics that were newly computed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) classFields: rewrite the data rows (header stays put) in the new
#    order, keeping the same four columns (Class Name / Field Name /
#    Field Modifier / Field Type).
# ---------------------------------------------------------------------
$fields = $wb.Worksheets.Item("classFields")

$fieldRows = @(
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "id", "private", "java.lang.Long"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "name", "private", "java.lang.String"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "availableItems", "private", "int"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "reservedItems", "private", "int"),
    @("com.zatribune.spring.ecommerce.stock.db.DevBootstrap", "repository", "private", "com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository"),
    @("com.zatribune.spring.ecommerce.stock.db.DevBootstrap", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "id", "private", "java.lang.Long"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "name", "private", "java.lang.String"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "reservedItems", "private", "int"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "availableItems", "private", "int"),
    @("com.zatribune.spring.ecommerce.stock.listener.OrderListener", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.stock.listener.OrderListener", "orderService", "private", "com.zatribune.spring.ecommerce.stock.service.OrderService"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "SOURCE", "private", "domain.OrderSource"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "template", "private", "org.springframework.kafka.core.KafkaTemplate"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "log", "private", "org.slf4j.Logger"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "repository", "private", "com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository")
)

for ($i = 0; $i -lt $fieldRows.Count; $i++) {
    $r = $i + 2
    $row = $fieldRows[$i]
    $fields.Cells.Item($r, 1).Value = $row[0]
    $fields.Cells.Item($r, 2).Value = $row[1]
    $fields.Cells.Item($r, 3).Value = $row[2]
    $fields.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2) New sheet: classNumberOfLines
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$classLines = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$classLines.Name = "classNumberOfLines"

$classLines.Range("A1").Value = "Class Name"
$classLines.Range("B1").Value = "Number of Lines"

$classLineRows = @(
    @("com.zatribune.spring.ecommerce.stock.StockApplicationTests", "5"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "31"),
    @("com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository", "3"),
    @("com.zatribune.spring.ecommerce.stock.db.DevBootstrap", "18"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "1"),
    @("com.zatribune.spring.ecommerce.stock.listener.OrderListener", "20"),
    @("com.zatribune.spring.ecommerce.stock.StockApplication", "6"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "45"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderService", "5")
)

for ($i = 0; $i -lt $classLineRows.Count; $i++) {
    $r = $i + 2
    $row = $classLineRows[$i]
    $classLines.Cells.Item($r, 1).Value = $row[0]
    $classLines.Cells.Item($r, 2).Value = "'" + $row[1]
}

# ---------------------------------------------------------------------
# 3) New sheet: methodNumberOfLines
# ---------------------------------------------------------------------
$methodLines = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $classLines)
$methodLines.Name = "methodNumberOfLines"

$methodLines.Range("A1").Value = "Class Name"
$methodLines.Range("B1").Value = "Method Signature"
$methodLines.Range("C1").Value = "Number of Lines"

$methodLineRows = @(
    @("com.zatribune.spring.ecommerce.stock.StockApplicationTests", "contextLoads()", "2"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "id(java.lang.Long)", "4"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "name(java.lang.String)", "4"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "availableItems(int)", "4"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "reservedItems(int)", "4"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "build()", "3"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product`$ProductBuilder", "toString()", "3"),
    @("com.zatribune.spring.ecommerce.stock.db.DevBootstrap", "run(java.lang.String[])", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "toString()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "builder()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "getId()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "getName()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "getAvailableItems()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "getReservedItems()", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "setId(java.lang.Long)", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "setName(java.lang.String)", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "setAvailableItems(int)", "1"),
    @("com.zatribune.spring.ecommerce.stock.db.entities.Product", "setReservedItems(int)", "1"),
    @("com.zatribune.spring.ecommerce.stock.listener.OrderListener", "onEvent(domain.Order)", "2"),
    @("com.zatribune.spring.ecommerce.stock.StockApplication", "main(java.lang.String[])", "3"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "reserve(domain.Order)", "3"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderServiceImpl", "confirm(domain.Order)", "3"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderService", "reserve(domain.Order)", "1"),
    @("com.zatribune.spring.ecommerce.stock.service.OrderService", "confirm(domain.Order)", "1")
)

for ($i = 0; $i -lt $methodLineRows.Count; $i++) {
    $r = $i + 2
    $row = $methodLineRows[$i]
    $methodLines.Cells.Item($r, 1).Value = $row[0]
    $methodLines.Cells.Item($r, 2).Value = $row[1]
    $methodLines.Cells.Item($r, 3).Value = "'" + $row[2]
}

$wb.Worksheets.Item("classMethods").Select()
